$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Intern shared strings in the same order as the target file: CONFIDENCE.T (idx 9) then CONFIDENCE.NORM (idx 10)
$ws.Range("A10").Value = "CONFIDENCE.T"
$ws.Range("A9").Value = "CONFIDENCE.NORM"

# Row 9: CONFIDENCE.NORM
$ws.Range("B9").Formula = "=_xlfn.CONFIDENCE.NORM(C9,D9,E9)"
$ws.Range("C9").Formula = "=2/15"
$ws.Range("D9").Value = 6.6
$ws.Range("E9").Value = 44

# Row 10: CONFIDENCE.T (label only; formula still references CONFIDENCE.NORM, matching source)
$ws.Range("B10").Formula = "=_xlfn.CONFIDENCE.NORM(C10,D10,E10)"
$ws.Range("C10").Formula = "=2/15"
$ws.Range("D10").Value = 6.6
$ws.Range("E10").Value = 44

# Column width change (A: 15.7109375 -> 19.28515625)
$ws.Columns.Item(1).ColumnWidth = 18.48

# Selection moves to A11 after the new rows
$ws.Range("A11").Select()
